$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4542502.5
$ws.Range("I33").Value = 7138091.5
$ws.Range("K33").Value = 7138091.5
$ws.Range("M33").Value = -7137862.5

$ws.Range("H41").Value = 286.36365
$ws.Range("I41").Value = 217.14285
$ws.Range("J41").Value = 407.5
$ws.Range("K41").Value = 217.14285
$ws.Range("L41").Value = 407.5
$ws.Range("M41").Value = 222.85715
$ws.Range("N41").Value = -1287.5

$ws.Range("H100").Value = 1694.6666
$ws.Range("I100").Value = 1640
$ws.Range("J100").Value = 1845
$ws.Range("K100").Value = 1640
$ws.Range("L100").Value = 1845
$ws.Range("M100").Value = -1099
$ws.Range("N100").Value = -2927

$ws.Range("H125").Value = 6334.8
$ws.Range("J125").Value = 6875.778
$ws.Range("L125").Value = 61882.002
$ws.Range("N125").Value = -66802.00200000001

$ws.Range("H138").Value = 4676.5537
$ws.Range("I138").Value = 2334.4167
$ws.Range("J138").Value = 5315.3184
$ws.Range("K138").Value = 7003.250100000001
$ws.Range("L138").Value = 15945.9552
$ws.Range("M138").Value = -1863.250100000001
$ws.Range("N138").Value = -26225.9552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2090.8696
$ws.Range("I45").Value = 2049.5454
$ws.Range("K45").Value = 2049.5454
$ws.Range("M45").Value = -1672.5454

$ws.Range("H61").Value = 5546.3145
$ws.Range("I61").Value = 4028.3572
$ws.Range("J61").Value = 11618.143
$ws.Range("K61").Value = 4028.3572
$ws.Range("L61").Value = 11618.143
$ws.Range("M61").Value = -3816.3572
$ws.Range("N61").Value = -12042.143

$ws.Range("H136").Value = 5546.3145
$ws.Range("I136").Value = 4028.3572
$ws.Range("J136").Value = 11618.143
$ws.Range("K136").Value = 12085.0716
$ws.Range("L136").Value = 34854.429
$ws.Range("M136").Value = -9535.071599999999
$ws.Range("N136").Value = -39954.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 460
$ws.Range("I22").Value = 460
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 460
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -287
$ws.Range("N22").ClearContents()

$ws.Range("H33").Value = 20031.75
$ws.Range("J33").Value = 20031.75
$ws.Range("L33").Value = 20031.75
$ws.Range("N33").Value = -20703.75

$ws.Range("H56").Value = 36110
$ws.Range("J56").Value = 36110
$ws.Range("L56").Value = 36110
$ws.Range("N56").Value = -37588

$ws.Range("H105").Value = 9249.75
$ws.Range("I105").Value = 14499.5
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 14499.5
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -12752.5
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 291.5
$ws.Range("I22").Value = 291.5
$ws.Range("K22").Value = 291.5
$ws.Range("M22").Value = 58.5

$ws.Range("H31").Value = 472774.4
$ws.Range("I31").Value = 5579.625
$ws.Range("J31").Value = 951948.5600000001
$ws.Range("K31").Value = 5579.625
$ws.Range("L31").Value = 951948.5600000001
$ws.Range("M31").Value = -5284.625
$ws.Range("N31").Value = -952538.5600000001

$ws.Range("H34").Value = 472774.4
$ws.Range("I34").Value = 5579.625
$ws.Range("J34").Value = 951948.5600000001
$ws.Range("K34").Value = 5579.625
$ws.Range("L34").Value = 951948.5600000001
$ws.Range("M34").Value = -5377.625
$ws.Range("N34").Value = -952352.5600000001

$ws.Range("H105").Value = 1023.3333
$ws.Range("I105").Value = 988
$ws.Range("K105").Value = 988
$ws.Range("M105").Value = 759

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38461748
$ws.Range("I12").Value = 142857380
$ws.Range("J12").Value = 202.21053
$ws.Range("K12").Value = 428572140
$ws.Range("L12").Value = 606.63159
$ws.Range("M12").Value = -428571967
$ws.Range("N12").Value = -952.63159

$ws.Range("H34").Value = 2639.111
$ws.Range("I34").Value = 332.66666
$ws.Range("J34").Value = 3100.4
$ws.Range("K34").Value = 997.9999799999999
$ws.Range("L34").Value = 9301.200000000001
$ws.Range("M34").Value = -913.9999799999999
$ws.Range("N34").Value = -9469.200000000001

$ws.Range("H39").Value = 5785.7144
$ws.Range("J39").Value = 5785.7144
$ws.Range("L39").Value = 17357.1432
$ws.Range("N39").Value = -17945.1432

$ws.Range("H55").Value = 7975

$ws.Range("H86").Value = 2850.5
$ws.Range("I86").Value = 2850.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 8551.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -7365.5
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 2850.5
$ws.Range("I89").Value = 2850.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25654.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19726.5
$ws.Range("N89").ClearContents()

$ws.Range("H97").Value = 158.2
$ws.Range("I97").Value = 133.33333
$ws.Range("J97").Value = 195.5
$ws.Range("K97").Value = 399.99999
$ws.Range("L97").Value = 586.5
$ws.Range("M97").Value = 96.00001000000003
$ws.Range("N97").Value = -1578.5

$ws.Range("H105").Value = 39422.223
$ws.Range("J105").Value = 39422.223
$ws.Range("L105").Value = 118266.669
$ws.Range("N105").Value = -123508.669

$ws.Range("H113").Value = 568.8099999999999
$ws.Range("I113").Value = 664.46
$ws.Range("J113").Value = 473.16
$ws.Range("K113").Value = 1993.38
$ws.Range("L113").Value = 1419.48
$ws.Range("M113").Value = 176.6199999999999
$ws.Range("N113").Value = -5759.48

$ws.Range("H137").Value = 21321.814
$ws.Range("I137").Value = 1158.1177
$ws.Range("J137").Value = 55600.1
$ws.Range("K137").Value = 3474.3531
$ws.Range("L137").Value = 166800.3
$ws.Range("M137").Value = 1625.6469
$ws.Range("N137").Value = -177000.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 10000000
$ws.Range("I24").Value = 10000000
$ws.Range("K24").Value = 10000000
$ws.Range("M24").Value = -9999827

$ws.Range("H39").Value = 38761
$ws.Range("J39").Value = 38761
$ws.Range("L39").Value = 38761
$ws.Range("N39").Value = -39825

$ws.Range("H113").Value = 2416.5833
$ws.Range("I113").Value = 2150
$ws.Range("J113").Value = 2949.75
$ws.Range("K113").Value = 2150
$ws.Range("L113").Value = 2949.75
$ws.Range("M113").Value = 20
$ws.Range("N113").Value = -7289.75

$ws.Range("H126").Value = 2819.7827
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 477.8
$ws.Range("I22").Value = 477.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 477.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -182.8
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 477.8
$ws.Range("I27").Value = 477.8
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 477.8
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -370.8
$ws.Range("N27").ClearContents()

$ws.Range("H53").Value = 23017
$ws.Range("J53").Value = 23017
$ws.Range("L53").Value = 23017
$ws.Range("N53").Value = -24053

$ws.Range("H61").Value = 27468.256
$ws.Range("I61").Value = 28328.7
$ws.Range("J61").Value = 24600.111
$ws.Range("K61").Value = 28328.7
$ws.Range("L61").Value = 24600.111
$ws.Range("M61").Value = -28126.7
$ws.Range("N61").Value = -25004.111

$ws.Range("H113").Value = 27468.256
$ws.Range("I113").Value = 28328.7
$ws.Range("J113").Value = 24600.111
$ws.Range("K113").Value = 28328.7
$ws.Range("L113").Value = 24600.111
$ws.Range("M113").Value = -26158.7
$ws.Range("N113").Value = -28940.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 10084
$ws.Range("J53").Value = 10084
$ws.Range("L53").Value = 10084
$ws.Range("N53").Value = -11298

$ws.Range("H61").Value = 6940.375
$ws.Range("I61").Value = 2714
$ws.Range("J61").Value = 13984.333
$ws.Range("K61").Value = 2714
$ws.Range("L61").Value = 13984.333
$ws.Range("M61").Value = -2422
$ws.Range("N61").Value = -14568.333

$ws.Range("H74").Value = 14344.9
$ws.Range("J74").Value = 14827.667
$ws.Range("L74").Value = 14827.667
$ws.Range("N74").Value = -16699.667

$ws.Range("H77").Value = 14344.9
$ws.Range("J77").Value = 14827.667
$ws.Range("L77").Value = 44483.001
$ws.Range("N77").Value = -53843.001
